# This script applies a cyclic rotation of data among rows 4, 5 and 6:
#   new row4 <- old row5
#   new row5 <- old row6
#   new row6 <- old row4
# Only the following columns actually differ between the three rows:
#   A, B, D, E, F, G, H, P, Q, R, AI
# All other columns are identical across rows 4/5/6, so only these columns
# need to be rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","D","E","F","G","H","P","Q","R","AI")

# Capture the original values for rows 4, 5 and 6 before overwriting anything.
# Note: Use .Value2 (not .Value) - in this runtime .Value round-trips through
# a Variant wrapper object that does not unwrap properly when re-assigned.
$orig = @{}
foreach ($r in 4,5,6) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# Apply the rotation: row4 <- old row5, row5 <- old row6, row6 <- old row4
$mapping = @{ 4 = 5; 5 = 6; 6 = 4 }

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $orig[$srcRow][$col]
    }
}
